$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill the previously-empty "Notes" column (E) cells with the literal
# string "nan" for the rows that did not already have a value.
$ws.Range("E2").Value = "nan"
$ws.Range("E3").Value = "nan"
$ws.Range("E4").Value = "nan"
$ws.Range("E5").Value = "nan"
$ws.Range("E6").Value = "nan"
$ws.Range("E7").Value = "nan"
$ws.Range("E8").Value = "nan"
$ws.Range("E12").Value = "nan"

# Match the author's final on-screen selection/scroll position.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("E40").Select()
